# Update cryptos list with latest prices and volume changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.144.20'
$ws.Range("E2").Value = '  +0.15%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.420.46'
$ws.Range("E3").Value = '  +0.00%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '553.62'

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.12'
$ws.Range("E6").Value = '  -0.47%  '

# Row 7
$ws.Range("E7").Value = '  +0.04%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.588'
$ws.Range("E8").Value = '  +2.13%  '

# Row 9
$ws.Range("E9").Value = '  -1.32%  '

# Row 10
$ws.Range("E10").Value = '  -0.79%  '

# Row 11
$ws.Range("E11").Value = '  -0.29%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.91'
$ws.Range("E13").Value = '  +0.07%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.855.21'
$ws.Range("E14").Value = '  +0.22%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '60.083.93'

# Row 16
$ws.Range("E16").Value = '  -0.30%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.420.13'
$ws.Range("E17").Value = '  +0.35%  '

# Row 18
$ws.Range("E18").Value = '  -0.45%  '

# Row 19
$ws.Range("E19").Value = '  +2.70%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '327.09'
$ws.Range("E20").Value = '  -1.45%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.75'
$ws.Range("E21").Value = '  +0.16%  '

# Row 22
$ws.Range("E22").Value = '  -0.06%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.18'
$ws.Range("E23").Value = '  +0.04%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.177'
$ws.Range("E24").Value = '  +4.04%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.64'
$ws.Range("E25").Value = '  +0.90%  '

# Row 26
$ws.Range("E26").Value = '  +0.23%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.41'
$ws.Range("E27").Value = '  +5.68%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0774'
$ws.Range("E28").Value = '  -1.14%  '

# Row 29
$ws.Range("E29").Value = '  -0.18%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '170.64'
$ws.Range("E30").Value = '  +0.68%  '

# Row 31
$ws.Range("E31").Value = '  -2.40%  '

# Row 32
$ws.Range("E32").Value = '  -3.27%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.06'
$ws.Range("E33").Value = '  +2.07%  '

# Row 34
$ws.Range("E34").Value = '  -0.65%  '

# Row 35
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.34'
$ws.Range("E35").Value = '  +2.87%  '

# Row 36
$ws.Range("B36").Value = 'USDe'
$ws.Range("C36").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  +0.02%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.22'
$ws.Range("E37").Value = '  +0.41%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '326.83'
$ws.Range("E39").Value = '  +4.46%  '

# Row 40
$ws.Range("E40").Value = '  -0.71%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '145.15'

# Row 42
$ws.Range("E42").Value = '  -0.82%  '

# Row 43
$ws.Range("B43").Value = 'Stellar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0963'
$ws.Range("E43").Value = '  +0.32%  '

# Row 44
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.88'
$ws.Range("E44").Value = '  +1.92%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0517'
$ws.Range("E45").Value = '  -0.68%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.575'
$ws.Range("E46").Value = '  +0.15%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0223'
$ws.Range("E47").Value = '  -0.99%  '

# Row 48
$ws.Range("E48").Value = '  -0.08%  '

# Row 49
$ws.Range("E49").Value = '  -1.40%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.66'
$ws.Range("E50").Value = '  -0.49%  '

# Row 51
$ws.Range("E51").Value = '  -0.63%  '
